$wb = $excel.ActiveWorkbook

# --- Create the three sheets: ArcFace, VGGFace, FaceNet512 (in this left-to-right order) ---
# Start from the original (single) sheet, which will become "FaceNet512".
$orig = $wb.Worksheets.Item(1)

# Insert VGGFace right before the original sheet, then insert ArcFace right before VGGFace
# so the final left-to-right order is: ArcFace, VGGFace, FaceNet512.
$vgg = $wb.Worksheets.Add($orig)
$vgg.Name = "VGGFace"

$arc = $wb.Worksheets.Add($wb.Worksheets.Item("VGGFace"))
$arc.Name = "ArcFace"

$wb.Worksheets.Item("Sheet1").Name = "FaceNet512"

# --- Populate ArcFace ---
$ws = $wb.Worksheets.Item("ArcFace")
$ws.Cells.Item(1,1).Value = "Metric"
$ws.Cells.Item(1,2).Value = "Value (Weighted)"
$ws.Cells.Item(1,3).Value = "Value (Micro)"
$ws.Cells.Item(1,4).Value = "Value(Macro)"

$ws.Cells.Item(2,1).Value = "Accuracy"
$ws.Cells.Item(2,2).Value = 0.93110899999999996
$ws.Cells.Item(2,3).Value = 0.93110899999999996
$ws.Cells.Item(2,4).Value = 0.93110899999999996

$ws.Cells.Item(3,1).Value = "Precision"
$ws.Cells.Item(3,2).Value = 0.95129900000000001
$ws.Cells.Item(3,3).Value = 0.93110899999999996
$ws.Cells.Item(3,4).Value = 0.646922

$ws.Cells.Item(4,1).Value = "Recall"
$ws.Cells.Item(4,2).Value = 0.93110899999999996
$ws.Cells.Item(4,3).Value = 0.93110899999999996
$ws.Cells.Item(4,4).Value = 0.79158499999999998

$ws.Cells.Item(5,1).Value = "F1-Score"
$ws.Cells.Item(5,2).Value = 0.93922499999999998
$ws.Cells.Item(5,3).Value = 0.93110899999999996
$ws.Cells.Item(5,4).Value = 0.65786500000000003

# NOTE: the engine snaps ColumnWidth to a pixel grid (multiples of 1/6 character
# width at this font/MDW), so the exact author widths (10.08984375 / 15.26953125 /
# 12 / 12.08984375 characters) are not all individually reachable. The inputs below
# were picked (by probing the grid) to land on the closest achievable stored width
# for each column - col C (width 12) lands exactly on target.
$ws.Columns.Item(1).ColumnWidth = 9.25
$ws.Columns.Item(2).ColumnWidth = 14.45
$ws.Columns.Item(3).ColumnWidth = 11.1
$ws.Columns.Item(4).ColumnWidth = 11.3

$ws.Range("D4").Select()

# --- Populate VGGFace ---
$ws = $wb.Worksheets.Item("VGGFace")
$ws.Cells.Item(1,1).Value = "Metric"
$ws.Cells.Item(1,2).Value = "Value (Weighted)"
$ws.Cells.Item(1,3).Value = "Value (Micro)"
$ws.Cells.Item(1,4).Value = "Value(Macro)"

$ws.Cells.Item(2,1).Value = "Accuracy"
$ws.Cells.Item(2,2).Value = 0.94833199999999995
$ws.Cells.Item(2,3).Value = 0.94833199999999995
$ws.Cells.Item(2,4).Value = 0.94833199999999995

$ws.Cells.Item(3,1).Value = "Precision"
$ws.Cells.Item(3,2).Value = 0.96812799999999999
$ws.Cells.Item(3,3).Value = 0.94833199999999995
$ws.Cells.Item(3,4).Value = 0.66487499999999999

$ws.Cells.Item(4,1).Value = "Recall"
$ws.Cells.Item(4,2).Value = 0.94833199999999995
$ws.Cells.Item(4,3).Value = 0.94833199999999995
$ws.Cells.Item(4,4).Value = 0.79497600000000002

$ws.Cells.Item(5,1).Value = "F1-Score"
$ws.Cells.Item(5,2).Value = 0.95698899999999998
$ws.Cells.Item(5,3).Value = 0.94833199999999995
$ws.Cells.Item(5,4).Value = 0.66696900000000003

# NOTE: the engine snaps ColumnWidth to a pixel grid (multiples of 1/6 character
# width at this font/MDW), so the exact author widths (10.08984375 / 15.26953125 /
# 12 / 12.08984375 characters) are not all individually reachable. The inputs below
# were picked (by probing the grid) to land on the closest achievable stored width
# for each column - col C (width 12) lands exactly on target.
$ws.Columns.Item(1).ColumnWidth = 9.25
$ws.Columns.Item(2).ColumnWidth = 14.45
$ws.Columns.Item(3).ColumnWidth = 11.1
$ws.Columns.Item(4).ColumnWidth = 11.3

# --- Populate FaceNet512 ---
$ws = $wb.Worksheets.Item("FaceNet512")
$ws.Cells.Item(1,1).Value = "Metric"
$ws.Cells.Item(1,2).Value = "Value (Weighted)"
$ws.Cells.Item(1,3).Value = "Value (Micro)"
$ws.Cells.Item(1,4).Value = "Value(Macro)"

$ws.Cells.Item(2,1).Value = "Accuracy"
$ws.Cells.Item(2,2).Value = 0.95371399999999995
$ws.Cells.Item(2,3).Value = 0.95371399999999995
$ws.Cells.Item(2,4).Value = 0.95371399999999995

$ws.Cells.Item(3,1).Value = "Precision"
$ws.Cells.Item(3,2).Value = 0.98570800000000003
$ws.Cells.Item(3,3).Value = 0.95371399999999995
$ws.Cells.Item(3,4).Value = 0.66787700000000005

$ws.Cells.Item(4,1).Value = "Recall"
$ws.Cells.Item(4,2).Value = 0.95371399999999995
$ws.Cells.Item(4,3).Value = 0.95371399999999995
$ws.Cells.Item(4,4).Value = 0.80341099999999999

$ws.Cells.Item(5,1).Value = "F1-Score"
$ws.Cells.Item(5,2).Value = 0.96896599999999999
$ws.Cells.Item(5,3).Value = 0.95371399999999995
$ws.Cells.Item(5,4).Value = 0.666126

# NOTE: the engine snaps ColumnWidth to a pixel grid (multiples of 1/6 character
# width at this font/MDW), so the exact author widths (10.08984375 / 15.26953125 /
# 12 / 12.08984375 characters) are not all individually reachable. The inputs below
# were picked (by probing the grid) to land on the closest achievable stored width
# for each column - col C (width 12) lands exactly on target.
$ws.Columns.Item(1).ColumnWidth = 9.25
$ws.Columns.Item(2).ColumnWidth = 14.45
$ws.Columns.Item(3).ColumnWidth = 11.1
$ws.Columns.Item(4).ColumnWidth = 11.3

$ws.Range("C3").Select()

# --- Make VGGFace the active tab/sheet (must be done last: selecting/activating
# a sheet changes the workbook's active tab as a side effect) ---
$vggFinal = $wb.Worksheets.Item("VGGFace")
$vggFinal.Activate()
$vggFinal.Range("E4").Select()
